# Commit: "SQL for inserting photos"
# Adds a new worksheet "20201026" (after "20201020") containing the same
# ID / Shop ID / SQL generator pattern as the other date sheets, for the
# shop "da04f82c-ffb0-11ea-ba65-065a10bcba76" with 21 photo rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Clear the previously-active sheet's (20201020 / sheet8) selection so
#    it no longer shows the old "mid-edit" cursor position, matching the
#    tidy-up that happens when the user moves on to a freshly added sheet.
# ---------------------------------------------------------------------
$prevActive = $wb.Worksheets.Item($wb.Worksheets.Count)
$prevActive.Activate()
$prevActive.Range("A1:C2").Select()

# ---------------------------------------------------------------------
# 2. Add the new sheet at the end of the tab strip and name it.
# ---------------------------------------------------------------------
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $prevActive)
$new.Name = "20201026"

# Match the (auto-fit) column B width used by all the other date sheets
# (best achievable value given this engine's pixel-width quantization).
$new.Columns.Item(2).ColumnWidth = 34.3

$shopId = "da04f82c-ffb0-11ea-ba65-065a10bcba76"

# ---------------------------------------------------------------------
# 3. Header row.
# ---------------------------------------------------------------------
$new.Cells.Item(1, 1).Value = "ID"
$new.Cells.Item(1, 2).Value = "Shop ID"
$new.Cells.Item(1, 3).Value = "SQL"
$new.Range("A1:C1").Font.Color = 0

# ---------------------------------------------------------------------
# 4. 21 data rows (ID 1..21), all for the same shop id, with the SQL
#    builder formula in column C.
# ---------------------------------------------------------------------
for ($i = 1; $i -le 21; $i++) {
    $row = $i + 1

    $new.Cells.Item($row, 1).Value = $i
    $new.Cells.Item($row, 2).Value = $shopId
    $new.Cells.Item($row, 3).Formula = '=_xlfn.CONCAT("INSERT INTO photos(restaurant_id, name, type) VALUES(UuidToBin(''", B' + $row + ', "''), LPAD(", A' + $row + ', ", 7, ''0''), ''dish''", ");")'

    # Banded-row formatting carried over from the other date sheets:
    # column A alternates on every even row, column B is only styled on
    # the first data row (row 2).
    if ($row % 2 -eq 0) {
        $new.Cells.Item($row, 1).Font.Color = 0
    }
    if ($row -eq 2) {
        $new.Cells.Item($row, 2).Font.Color = 0
    }
}

# ---------------------------------------------------------------------
# 5. Final selection / active-cell state for the new sheet and make it
#    the active tab.
# ---------------------------------------------------------------------
$new.Activate()
$new.Range("B11").Select()
